$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" numeric-looking string (e.g. "41.36")
# must be forced to stay text (matching the source inlineStr cells) by
# temporarily switching the cell to Text format, then restoring the default
# "Normal" style so no stray style/number-format diff is introduced.

$ws.Range("D2").Value = "35.237.76"
$ws.Range("D3").Value = "1.894.31"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("E6").Value = "  +5.51%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.36"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.348"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0711"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0994"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "2.168.42"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.05%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.904.40"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.690"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "35.251.66"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "0.0₃0815"
$ws.Range("E20").Value = "  +2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "240.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +29.34%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.943"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0561"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0210"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0648"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "89.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").Value = "1.336.57"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("B46").Value = "MultiversX"
$ws.Range("C46").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +38.27%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "2.079.15"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.03%  "
